$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit cyclically rotates the data of rows 6-9 (the actual
# observation values), while keeping several shared/unchanged columns
# (C, D, I, J, K, L, N, O, S, T, U, V, W, X, Y, Z, AA, AB, AD-AY, etc.) as
# they already are (they happen to be identical across these rows, or
# are untouched by the source diff).
#
# Effective mapping of row *content*:
#   new row 6  <= old row 9
#   new row 7  <= old row 8
#   new row 8  <= old row 6
#   new row 9  <= old row 7
#
# Rather than rewriting every single cell, only the cells that actually
# differ per the diff are updated explicitly below.

# --- Row 6 (receives old row 9 species/finding data) ---
$ws.Range("A6").Value2 = 131017110
$ws.Range("B6").Value2 = 57884
$ws.Range("E6").Value2 = 100109
$ws.Range("F6").Value2 = "Tretåig hackspett"
$ws.Range("G6").Value2 = "Picoides tridactylus"
$ws.Range("H6").Value2 = "(Linnaeus, 1758)"
$ws.Range("M6").Value2 = "färska spår"
$ws.Range("Q6").Value2 = 477185
$ws.Range("R6").Value2 = 6789174

# --- Row 7 (receives old row 8 coordinates/id; loses its comment) ---
$ws.Range("A7").Value2 = 131016886
$ws.Range("Q7").Value2 = 477116
$ws.Range("R7").Value2 = 6789167
$ws.Range("AC7").ClearContents()

# --- Row 8 (receives old row 6 coordinates/id) ---
$ws.Range("A8").Value2 = 131016558
$ws.Range("Q8").Value2 = 477128
$ws.Range("R8").Value2 = 6789106

# --- Row 9 (receives old row 7 species/finding data plus its comment) ---
$ws.Range("A9").Value2 = 131017563
$ws.Range("B9").Value2 = 79243
$ws.Range("E9").Value2 = 6425
$ws.Range("F9").Value2 = "Garnlav"
$ws.Range("G9").Value2 = "Alectoria sarmentosa"
$ws.Range("H9").Value2 = "(Ach.) Ach."
$ws.Range("M9").ClearContents()
$ws.Range("Q9").Value2 = 477226
$ws.Range("R9").Value2 = 6789084
$ws.Range("AC9").Value2 = "Rikligt i området"
